$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 9 so "US7" (row 8) can be split into
#     "US7RejectOrder" (row 8) and a new "US7VerifyOrder" (row 9),
#     pushing the old "US8" row down to row 10 where it becomes
#     "US8AmendOrder".
$ws.Rows.Item(9).Insert()

$ws.Range("A8").Value() = "US7RejectOrder"
$ws.Range("A9").Value() = "US7VerifyOrder"
$ws.Range("A10").Value() = "US8AmendOrder"

# --- The row that held "US11" has, after the insertion above, shifted down
#     to row 13; turn it into the new "US11TimeOut" test case in place.
$ws.Range("A13").Value() = "US11TimeOut"
$ws.Range("B13").Value() = "Testing this can be done. Commented out because of long wait time."
$ws.Range("C13").Value() = "PASS"

# --- Extend the conditional formatting range on column C by one row
#     (C2:C627 -> C2:C628) to account for the inserted row shifting
#     the sheet's effective extent.
$cf = $ws.Range("C2:C627").FormatConditions
for ($i = 1; $i -le $cf.Count; $i++) {
    $cf.Item($i).ModifyAppliesToRange($ws.Range("C2:C628"))
}

# --- Update the active selection to match the saved view state.
$ws.Activate() | Out-Null
$ws.Range("B13").Select() | Out-Null
